$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.831.39'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '3.452.17'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'580.93"
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").Value = "'150.67"
$ws.Range("E6").Value = '  +2.30%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +1.26%  '
$ws.Range("D9").Value = "'8.09"
$ws.Range("E9").Value = '  +6.19%  '
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("D11").Value = "'0.418"
$ws.Range("E11").Value = '  +4.55%  '
$ws.Range("D12").Value = '4.038.47'
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = "'28.35"
$ws.Range("E14").Value = '  -4.44%  '
$ws.Range("D15").Value = '3.440.16'
$ws.Range("E15").Value = '  -1.09%  '
$ws.Range("D16").Value = "'0.0000174"
$ws.Range("E16").Value = '  +1.29%  '
$ws.Range("D17").Value = '62.773.52'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = "'6.45"
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("D19").Value = "'14.61"
$ws.Range("E19").Value = '  +1.70%  '
$ws.Range("D20").Value = "'9.01"
$ws.Range("E20").Value = '  -2.36%  '
$ws.Range("D21").Value = "'388.31"
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("D22").Value = "'0.568"
$ws.Range("E22").Value = '  +1.13%  '
$ws.Range("D23").Value = "'75.23"
$ws.Range("E23").Value = '  +0.49%  '
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = "'0.0000116"
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '3.586.88'
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("D27").Value = "'0.184"
$ws.Range("E27").Value = '  +2.63%  '
$ws.Range("D28").Value = "'7.77"
$ws.Range("E28").Value = '  +1.80%  '
$ws.Range("D30").Value = "'8.04"
$ws.Range("E30").Value = '  -1.44%  '
$ws.Range("D31").Value = "'2.13"
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").Value = "'1.36"
$ws.Range("E33").Value = '  -2.72%  '
$ws.Range("D34").Value = "'23.30"
$ws.Range("E34").Value = '  -1.87%  '
$ws.Range("D35").Value = "'5.45"
$ws.Range("E35").Value = '  +2.43%  '
$ws.Range("D36").Value = "'1.65"
$ws.Range("E36").Value = '  +3.52%  '
$ws.Range("D37").Value = "'31.86"
$ws.Range("E37").Value = '  +1.18%  '
$ws.Range("D38").Value = "'6.97"
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("D39").Value = "'169.03"
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("D40").Value = '3.481.70'
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("D41").Value = "'0.0785"
$ws.Range("E41").Value = '  +2.27%  '
$ws.Range("D42").Value = "'42.86"
$ws.Range("E42").Value = '  +1.40%  '
$ws.Range("D43").Value = "'0.784"
$ws.Range("E43").Value = '  -2.08%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = "'1.70"
$ws.Range("E44").Value = '  -1.24%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = "'4.40"
$ws.Range("E45").Value = '  -1.91%  '
$ws.Range("D46").Value = "'1.18"
$ws.Range("E46").Value = '  -1.34%  '
$ws.Range("D47").Value = '2.565.01'
$ws.Range("E47").Value = '  -1.53%  '
$ws.Range("D48").Value = "'6.92"
$ws.Range("E48").Value = '  +2.69%  '
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("D50").Value = "'22.91"
$ws.Range("E50").Value = '  -2.08%  '
$ws.Range("E51").Value = '  -0.02%  '
